# vault backup: 2024-11-20 08:04:06
# Updates the grad-roles tracker:
#  - re-colour a few rows whose deadlines have now passed / are rolling
#  - add a new "Postgraduate Intern" (Bank of England) row
#  - freeze the header row and leave the selection near the bottom of the list

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Row 2 now reads like the "closed / rolling" rows (grey fill, style group used by row 3) ---
$ws.Range("A3:B3").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# --- Rows 5 and 7 move from the green "open" fill to the red "closed" fill (style group used by row 12) ---
$ws.Range("A12:B12").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D12").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E12").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Range("A12:B12").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("D12").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E12").Copy()
$ws.Range("E7").PasteSpecial(-4122)

# --- New row 20: Postgraduate Intern @ Bank of England ---
$ws.Range("D20").Value = "https://eoff.fa.em1.ukg.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_3001/job/670"
$ws.Range("A20").Value = "Postgraduate Intern"
$ws.Range("B20").Value = "Bank of England"
$ws.Range("C20").Value = 45809
$ws.Range("E20").Value = "London"
$ws.Range("G20").Value = "No"
$ws.Range("F20").Formula = "=C20-TODAY()"
$ws.Range("F19").Copy()
$ws.Range("F20").PasteSpecial(-4122)

# --- Freeze the header row and leave the selection on A21, like the saved view ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A21").Select()
